$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert new introductory row (row 2: 0 - Introduction to the tutorial)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Introduction to the tutorial"

# Re-order the tutorials list (rows 6-13) and flag some with "*" in column C
$ws.Range("B6").Value = "Mahalanobis"
$ws.Range("C6").Value = "*"

$ws.Range("B7").Value = "e.ellipse"
$ws.Range("C7").Value = "*"

$ws.Range("B8").Value = "in.el"
$ws.Range("C8").Value = "*"

$ws.Range("B9").Value = "fitNiche"

$ws.Range("B11").Value = "nicheG"

$ws.Range("C12").Value = "*"

$ws.Range("C13").Value = "*"

$ws.Range("B13").Select()
